# Various updates made by 10/21
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow column A like the reference edit did.
$ws.Columns("A").ColumnWidth = 4.63

# New "Consumption kWh" column (D), mirroring the style of the existing
# header in C1 and adding a value per month row.
$ws.Range("D1").Value = "Consumption kWh"
$ws.Range("D1").Style = $ws.Range("C1").Style

$ws.Range("D2").Value = 130.23
$ws.Range("D3").Value = 101.34
$ws.Range("D4").Value = 145.6
$ws.Range("D5").Value = 199.29
$ws.Range("D6").Value = 178.62
$ws.Range("D7").Value = 169.69
$ws.Range("D8").Value = 220.41
$ws.Range("D9").Value = 178.37
$ws.Range("D10").Value = 130.65
$ws.Range("D11").Value = 125.62
$ws.Range("D12").Value = 120.54
$ws.Range("D13").Value = 135.34

$dRange = $ws.Range("D2:D13")
$dRange.Font.ThemeColor = 1

# A few consumption (C column) figures were revised.
$ws.Range("C2").Value = 1490.0
$ws.Range("C3").Value = 1236.55
$ws.Range("C13").Value = 1713.23

# Restrict entries in the Consumption kWh column to decimals in range.
$cRange = $ws.Range("C2:C13")
$cRange.Validation.Add(2, 1, 1, 0.1, 10000.0)
$cRange.Validation.InputTitle = ""
$cRange.Validation.ErrorTitle = ""
$cRange.Validation.InputMessage = "Enter a number between .1 and 10000"
$cRange.Validation.ErrorMessage = ""
$cRange.Validation.ShowInput = $true
$cRange.Validation.ShowError = $true
